$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to reflect repulled data / mean calculation
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -8
$ws.Range("F8").Value = -7
$ws.Range("F9").Value = -6
